$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Penawaran" as the last sheet (after "Seleksi")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Penawaran"

# Populate header row
$newSheet.Range("A1").Value = "Kode Pekerjaan"
$newSheet.Range("B1").Value = "Benefit"

# Populate data row
$newSheet.Range("A2").Value = "L001"
$newSheet.Range("B2").Value = "Asuransi Kecelakaan Kerja"

# Restore the previously active sheet (Wawancara) as the selected tab,
# since adding a worksheet activates the new one by default
$wb.Worksheets.Item("Wawancara").Activate()
